$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the hidden "_GoBack" bookmark from the (now empty) paragraph right
#    after the team-members table to the paragraph that starts the
#    "* Asset:" line further down the document.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$target = $d.Content
$target.Find.Execute("* Asset:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(1)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null

# ---------------------------------------------------------------------------
# 2) "* The Library controls its assets (book, magazine, DVD,...)"
#    -> split so "DVD" moves from the end of the first run to the start of
#       the run that holds the trailing ",..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("magazine, DVD", $true, $false, $false, $false, $false, $true, 1, $false, "magazine, ", 2) | Out-Null
$d.Content.Find.Execute(",…", $true, $false, $false, $false, $false, $true, 1, $false, "DVD,…", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) ". Type" -> ". " + "Asset" + "Type" (the new "AssetType" word is wrapped
#    in spell-check proofErr markers, split across two runs).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(". Type", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$xmlFrag = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:rsidR='002F4DDF' w:rsidRDefault='002F4DDF' w:rsidP='002F4DDF'><w:pPr><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:color w:val='616161'/><w:sz w:val='27'/><w:szCs w:val='27'/><w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:color w:val='616161'/><w:sz w:val='27'/><w:szCs w:val='27'/><w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/></w:rPr><w:tab/><w:t xml:space='preserve'>. </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:color w:val='616161'/><w:sz w:val='27'/><w:szCs w:val='27'/><w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/></w:rPr><w:t>Asset</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:color w:val='616161'/><w:sz w:val='27'/><w:szCs w:val='27'/><w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/></w:rPr><w:t>Type</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$rng.InsertXML($xmlFrag)
